# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values on the active worksheet to the newly
# computed strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 7
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 2
    26 = 2
    27 = 1
    28 = 0
    29 = 6
    30 = 3
    31 = 6
    33 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
